# Re-brand the StructureDefinition spreadsheet from the old Alvearie/IBM
# identity to LinuxForHealth, bump the profile version, refresh the publish
# date, and fix up the two "Elements" table cells that carried the old
# canonical URL / a stale duplicated constraint.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: Property/Value summary table -----------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-path"

# Version: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet: per-element definition table ------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 ("Extension" root row): the ele-1/ext-1 Constraint(s) text was a
# stale duplicate of the one on the "Extension.extension" row (row 4) -
# clear it here.
$elements.Range("AI2").Value = ""

# Row 5 ("Extension.url"): the "Fixed Value" cell holds the profile's own
# canonical URL (Extension.url is fixed to the defining StructureDefinition
# URL) - update it to match the new host too.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-path"
